$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; existing rows 15-30 shift down to 16-31.
$ws.Rows("15").Insert()

# Populate the newly inserted row 15 with the new observation.
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Femacal de La Calera"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = "2021-10-15"
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 100112022
$ws.Range("G15").Value = "Arveja Verde"
$ws.Range("H15").Value = "Perfection"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 71
$ws.Range("K15").Value = 29000
$ws.Range("L15").Value = 30000
$ws.Range("M15").Value = 29507
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 1180
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
